# Add a new "aug18" worksheet at the end of the workbook (after "jul18"),
# matching the pattern used by the existing monthly score-history sheets.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "aug18"

# Header / first round block (round date 8/5/2018 -> serial 43317)
$ws.Range("A1").Value = 43317
$ws.Range("A1").NumberFormat = "d-mmm-yy"
$ws.Range("B1").Value = "Score"
$ws.Range("C1").Value = "Fairway"
$ws.Range("D1").Value = "GIR"
$ws.Range("E1").Value = "Putts"
$ws.Range("F1").Value = "Comment"

$ws.Range("A2").Value = "Hole 1"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = "R"
$ws.Range("E2").Value = 3

$ws.Range("A3").Value = "Hole 2"
$ws.Range("B3").Value = 4
$ws.Range("E3").Value = 2

$ws.Range("A4").Value = "Hole 3"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "S"
$ws.Range("E4").Value = 2

$ws.Range("A5").Value = "Hole 4"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "S"
$ws.Range("E5").Value = 2

$ws.Range("A6").Value = "Hole 5"
$ws.Range("B6").Value = 5
$ws.Range("E6").Value = 3

$ws.Range("A7").Value = "Hole 6"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "S"
$ws.Range("E7").Value = 2

$ws.Range("A8").Value = "Hole 7"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = "R"
$ws.Range("E8").Value = 1

$ws.Range("A9").Value = "Hole 8"
$ws.Range("B9").Value = 6
$ws.Range("E9").Value = 2

$ws.Range("A10").Value = "Hole 9"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = "R"
$ws.Range("E10").Value = 1

$ws.Range("A11").Value = "Hole 10"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "R"
$ws.Range("E11").Value = 1

$ws.Range("A12").Value = "Hole 11"
$ws.Range("B12").Value = 3
$ws.Range("E12").Value = 1

$ws.Range("A13").Value = "Hole 12"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "S"
$ws.Range("E13").Value = 2

$ws.Range("A14").Value = "Hole 13"
$ws.Range("B14").Value = 6
$ws.Range("C14").Value = "L"
$ws.Range("E14").Value = 2

$ws.Range("A15").Value = "Hole 14"
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "L"
$ws.Range("E15").Value = 2

$ws.Range("A16").Value = "Hole 15"
$ws.Range("B16").Value = 3
$ws.Range("E16").Value = 1

$ws.Range("A17").Value = "Hole 16"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = "S"
$ws.Range("E17").Value = 2

$ws.Range("A18").Value = "Hole 17"
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = "R"
$ws.Range("E18").Value = 2

$ws.Range("A19").Value = "Hole 18"
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = "S"
$ws.Range("E19").Value = 2

$ws.Range("B20").Formula = "=SUM(B2:B19)"
$ws.Range("E20").Formula = "=SUM(E2:E19)"

# Blank, date-formatted separator rows awaiting future rounds (same
# pattern as the other monthly sheets, e.g. jul18's rows 22/43/64).
$ws.Range("A22").NumberFormat = "d-mmm-yy"
$ws.Range("A43").NumberFormat = "d-mmm-yy"
$ws.Range("A64").NumberFormat = "d-mmm-yy"

# Column A best-fits to the "Hole NN" labels, same as the other sheets.
$ws.Columns.Item(1).AutoFit()

# Leave the cursor on C2, matching the saved selection state.
$ws.Range("C2").Select()
